# Merge the run-per-word title text into a single run by replacing the
# whole range through a Characters() sub-range (preserves the paragraph's
# pPr and the run's (empty) rPr).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, $titleRange.Length).Text = "Here is a single header"

# Likewise merge the speaker-notes run-per-word text into a single run.
$notesRange = $s.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesRange.Text = "and here are some notes"
